$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.785.61'
$ws.Range('E2').Value = '  -1.14%  '
$ws.Range('D3').Value = '2.322.99'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.75'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.10'
$ws.Range('E6').Value = '  -3.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.500'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '33.97'
$ws.Range('E10').Value = '  -4.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0780'
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '18.58'
$ws.Range('E12').Value = '  -4.89%  '
$ws.Range('E13').Value = '  +1.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.67'
$ws.Range('E14').Value = '  -3.93%  '
$ws.Range('D15').Value = '2.686.24'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = '2.311.51'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.788'
$ws.Range('E17').Value = '  -0.16%  '
$ws.Range('D18').Value = '42.735.05'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.08'
$ws.Range('E19').Value = '  -4.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.16'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').Value = '0.0₃0885'
$ws.Range('E21').Value = '  -1.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.84'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.96'
$ws.Range('E23').Value = '  -0.92%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.22'
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.41'
$ws.Range('E26').Value = '  -1.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.39'
$ws.Range('E27').Value = '  -2.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.35'
$ws.Range('E28').Value = '  +13.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.09'
$ws.Range('E29').Value = '  -0.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.19'
$ws.Range('E30').Value = '  -6.50%  '
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0755'
$ws.Range('E32').Value = '  +8.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.97'
$ws.Range('E33').Value = '  -0.93%  '
$ws.Range('B34').Value = 'Celestia'
$ws.Range('C34').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.15'
$ws.Range('E34').Value = '  -4.50%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '127.54'
$ws.Range('E35').Value = '  -22.68%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.35'
$ws.Range('E36').Value = '  -4.26%  '
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.32'
$ws.Range('E37').Value = '  -1.00%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.81'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('E39').Value = '  -1.04%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.01'
$ws.Range('E40').Value = '  +20.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.71'
$ws.Range('E41').Value = '  -3.15%  '
$ws.Range('E42').Value = '  -1.55%  '
$ws.Range('D43').Value = '1.922.69'
$ws.Range('E43').Value = '  -3.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0281'
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.05'
$ws.Range('E45').Value = '  -6.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.08'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.71'
$ws.Range('E47').Value = '  -2.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.87'
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('D49').Value = '2.552.83'
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '52.46'
$ws.Range('E50').Value = '  -2.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.52'
$ws.Range('E51').Value = '  -0.93%  '
